# Update copyright license date
# ---------------------------------------------------------------
# The xlsx was re-saved by Excel: the "6x6" sheet's view scrolled/
# re-selected a different range, and the "rand" sheet's volatile
# RAND() formulas recalculated to fresh values (handled automatically
# by the engine's load/recalc pipeline - every formula is recomputed
# when the workbook is loaded and again on save).
#
# Reproduce the view-state change that Excel recorded for the "6x6"
# sheet: the visible window scrolled so row 19 is at the top, and the
# selection moved from A62 to T27:W31.

$wb = $excel.ActiveWorkbook
$ws6x6 = $wb.Worksheets.Item("6x6")
$ws6x6.Activate()

# Scroll the window so A19 becomes the top-left visible cell.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1

# Move the selection to T27:W31 (active cell T27).
$ws6x6.Range("T27:W31").Select()

# Force a full recalculation so the volatile RAND() formulas on the
# "rand" sheet pick up fresh cached values (mirrors what Excel does
# every time the workbook is opened/edited and saved).
$excel.CalculateFull()
